# Update the "Date" metadata value and set the "Case Sensitive" value to the
# literal text "true" on the Metadata sheet of the workbook.

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

# Row 8: Date property -> update the timestamp value
$metadata.Range("B8").Value = "2023-09-01T08:48:57+00:00"

# Row 14: Case Sensitive property -> set its value to the text "true".
# A bare Value assignment of "true" gets auto-coerced to the Boolean TRUE by
# Excel's smart-typing, so instead enter it as a formula that evaluates to
# the text string "true", then convert that formula to a static value via
# copy / paste-special-values. This preserves the cell's existing style
# (no new text/quote-prefix style gets introduced) while keeping the stored
# cell type as a plain shared string.
$caseSensitiveCell = $metadata.Range("B14")
$caseSensitiveCell.Formula = "=""true"""
$caseSensitiveCell.Copy()
$caseSensitiveCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
